# Add specimen field adjacent_diseases
#
# Inserts three new columns (TEXT / ONTOLOGY / ONTOLOGY LABEL) for the new
# "adjacent_diseases" specimen_from_organism field, directly after the
# existing "KNOWN DISEASE(S)..." columns (U) and before "AUTOLYSIS SCORE"
# (old V), on the "Specimen from organism" worksheet. Also bumps the
# specimen_from_organism schema version referenced on the "Schemas" sheet
# from 10.5.0 to 10.6.0.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Specimen from organism")

# Insert 3 blank columns before the old column V ("AUTOLYSIS SCORE").
# The new columns inherit formatting from the column immediately to their
# left (U), which already carries the right header/description/example
# styles used throughout this template.
$ws.Range("V1:X1").EntireColumn().Insert()

# Row 4 (machine field-path values) is hidden on this sheet; temporarily
# unhide it while writing values so Excel doesn't stamp an autofit height
# on the hidden row, then restore it to hidden afterwards.
$ws.Rows.Item(4).Hidden() = $false

# Column V: adjacent_diseases.text
$ws.Range("V1").Value() = "ADJACENT DISEASE(S)"
$ws.Range("V2").Value() = "Short description of the disease(s) adjacent to the specimen's collection site (e.g. breast cancer)."
$ws.Range("V3").Value() = "If a healthy specimen is sampled from a site adjacent to diseased tissue, enter that tissue's disease here. If no diseased tissue is adjacent to the specimen, leave blank. For example: type 2 diabetes mellitus; normal; hepatic steatosis"
$ws.Range("V4").Value() = "specimen_from_organism.adjacent_diseases.text"

# Column W: adjacent_diseases.ontology
$ws.Range("W1").Value() = "ADJACENT DISEASE(S)) ONTOLOGY ID"
$ws.Range("W2").Value() = "An ontology term identifier in the form prefix:accession."
$ws.Range("W3").Value() = " For example: MONDO:0005148; PATO:0000461; HP:0001397"
$ws.Range("W4").Value() = "specimen_from_organism.adjacent_diseases.ontology"

# Column X: adjacent_diseases.ontology_label
$ws.Range("X1").Value() = "ADJACENT DISEASE(S) ONTOLOGY LABEL ONTOLOGY ID"
$ws.Range("X2").Value() = "The preferred label for the ontology term referred to in the ontology field. This may differ from the user-supplied value in the text field."
$ws.Range("X3").Value() = " For example: type 2 diabetes mellitus; normal; Hepatic steatosis"
$ws.Range("X4").Value() = "specimen_from_organism.adjacent_diseases.ontology_label"

$ws.Rows.Item(4).Hidden() = $true

# Bump the specimen_from_organism schema version referenced on the
# "Schemas" lookup sheet (10.5.0 -> 10.6.0) now that the new field exists.
$schemas = $wb.Worksheets.Item("Schemas")
$schemas.Range("A22").Value() = "https://schema.humancellatlas.org/type/biomaterial/10.6.0/specimen_from_organism"

# A handful of other sheets already had their machine-readable row 4
# hidden; keep the remaining "visible row 4" sheets consistent with that
# convention now that the template has been touched.
$sheetsToHideRow4 = @(
    "Project",
    "Project - Publications",
    "Cell suspension",
    "Sequence file",
    "Treatment protocol",
    "Sequencing protocol",
    "Analysis file",
    "Analysis protocol"
)
foreach ($name in $sheetsToHideRow4) {
    $s = $wb.Worksheets.Item($name)
    $s.Rows.Item(4).Hidden() = $true
}
